# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Day value (serial date) and all hourly price values for row 2
$ws.Range("A2").Value = 45897

$ws.Range("B2").Value = 97.15000000000001
$ws.Range("C2").Value = 90.14
$ws.Range("D2").Value = 86.06999999999999
$ws.Range("E2").Value = 73.25
$ws.Range("F2").Value = 69.40000000000001
$ws.Range("G2").Value = 76.42
$ws.Range("H2").Value = 83.38
$ws.Range("I2").Value = 90.38
$ws.Range("J2").Value = 90.38
$ws.Range("K2").Value = 69.40000000000001
$ws.Range("L2").Value = 17.35
$ws.Range("M2").Value = 9.65
$ws.Range("N2").Value = 5.64
$ws.Range("O2").Value = 4.31
$ws.Range("P2").Value = 0.65
$ws.Range("Q2").Value = 0.01
$ws.Range("R2").Value = 0.65
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 5.11
$ws.Range("U2").Value = 65.28
$ws.Range("V2").Value = 90
$ws.Range("W2").Value = 97.08
$ws.Range("X2").Value = 93.77
$ws.Range("Y2").Value = 88.09
$ws.Range("Z2").Value = 54.38

# Slot_4h (AA/AB) — max value unchanged, price updated
$ws.Range("AB2").Value = 92.23999999999999

# Slot_2h_frist (AC/AD) — label & price updated
$ws.Range("AC2").Value = "0h-2h"
$ws.Range("AD2").Value = 93.65000000000001

# Slot_2h_second (AE/AF) — label & price updated
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 93.54000000000001
